$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date column G (week of 2021-04-28 / serial 44314), matching the
#     existing date-header style used by D4:F4 (style "3") ---
$ws.Range("G4").Value = 44314
$ws.Range("F4").Copy()
$ws.Range("G4").PasteSpecial(-4122)   # xlPasteFormats - reuse style 3, no new style

$checkmark = [char]0x2713

# --- G5: first attendance checkmark in the new column. Starting from the
#     existing centered/size-14 "D/E/F" style (style 4) and only touching
#     the font color creates exactly one new font + one new cell style
#     (matching the diff's new font/cellXfs entries), reusing the shared
#     string already used by D/E/F (index 12, the checkmark glyph). ---
$ws.Range("F5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$g5 = $ws.Range("G5")
$g5.Value = $checkmark
$g5.Font.Color = 0

# --- G6-G9, G11, G12, G14, G15: reuse the newly minted style from G5 ---
$g5.Copy()
foreach ($r in 6,7,8,9,11,12,14,15) {
    $cell = $ws.Range("G$r")
    $cell.PasteSpecial(-4122)
    $cell.Value = $checkmark
}

# --- G10, G13: keep the plain existing "D/E/F" style (style 4) ---
foreach ($r in 10,13) {
    $src = $ws.Range("F$r")
    $src.Copy()
    $cell = $ws.Range("G$r")
    $cell.PasteSpecial(-4122)
    $cell.Value = $checkmark
}

# --- Selection moved to I14 after editing ---
[void]$ws.Range("I14").Select()
